$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B column: updated token labels (only rows whose token changed per diff)
$ws.Range("B2").Value = "<zero>"
$ws.Range("B3").Value = "<part>"
$ws.Range("B5").Value = "<water>"
$ws.Range("B6").Value = "<can>"
$ws.Range("B7").Value = "<papa>"
$ws.Range("B9").Value = "<number>"
$ws.Range("B11").Value = "<word>"
$ws.Range("B13").Value = "<backspace>"
$ws.Range("B15").Value = "<can>"

# C column: updated counts for every data row
$ws.Range("C2").Value = 20
$ws.Range("C3").Value = 24
$ws.Range("C4").Value = 26
$ws.Range("C5").Value = 36
$ws.Range("C6").Value = 28
$ws.Range("C7").Value = 34
$ws.Range("C8").Value = 34
$ws.Range("C9").Value = 21
$ws.Range("C10").Value = 30
$ws.Range("C11").Value = 31
$ws.Range("C12").Value = 29
$ws.Range("C13").Value = 25
$ws.Range("C14").Value = 32
$ws.Range("C15").Value = 9
